$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14: date changed (28.05.2020 -> 28.05.2021) ---
$ws.Range("A14").Value = 44344

# --- Row 15: new date, new hours entry, recompute, new activity text ---
$ws.Range("A15").Value = 44351
$ws.Range("B15").Value = 1
$ws.Range("D15").Value = "Internes Meeting"

# --- Row 16: new date, new hours entry, new activity text ---
$ws.Range("A16").Value = 44351
$ws.Range("B16").Value = 2
$ws.Range("D16").Value = "Debugging des HTTP-Servers im RedpitayaStub (Resultat: HTTP-Server wurde entfernt, nur WebSocket bleibt="

# --- Row 17: new date, new hours entry, new activity text ---
$ws.Range("A17").Value = 44353
$ws.Range("B17").Value = 6
$ws.Range("D17").Value = "Implementierung AcquirerOptions: Einheitenlogik ist vollständig, Inputvalidierung noch zu tun"

# --- Move the active selection to D18 (as last edited/selected cell) ---
$ws.Range("D18").Select() | Out-Null
